$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same set applied to every data row, B:Q)
$values = @(
    0.445980542904741,
    -0.05067851974480586,
    0.6185305930822008,
    -0.1467729861727587,
    0.4898227215065937,
    0.3288898821261522,
    0.6237281562693088,
    0.4168592253980262,
    0.2588347001794138,
    0.33784696278872,
    0.279043433096329,
    0.5734892170966707,
    0.05025235926527027,
    0.5979038403548141,
    22.22406457805329,
    34.41282282673529
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
